$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Metadata")

# Add the required "Experimental" boolean value (stored as text "true") in column B, row 7.
# Assigning the literal text "true" directly would be auto-interpreted as a native Excel
# boolean, so we use the apostrophe text-prefix to force a text value, then reapply the
# formatting of a neighboring "Value" cell so the style index is not altered.
$ws.Range("B7").Value = "'true"
$ws.Range("A6").Copy()
$ws.Range("B7").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Update the "Date" value in column B, row 8 to the new timestamp
$ws.Range("B8").Value = "2023-02-01T09:05:11-06:00"
